# Rename the column headers in row 1 so that the "_old" / "_new" suffixes
# are replaced by the respective format-version suffixes ("_FV2310" /
# "_FV2404"), add an Excel Table (ListObject) over the used range, and
# freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base header names (without suffix), in column order A..J (1..10)
$baseHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J (1-10): "<name>_old" -> "<name>_FV2310"
for ($i = 0; $i -lt $baseHeaders.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($baseHeaders[$i])_FV2310"
}

# Column K (11) "diff" is left untouched.

# Columns L-U (12-21): "<name>_new" -> "<name>_FV2404"
for ($i = 0; $i -lt $baseHeaders.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($baseHeaders[$i])_FV2404"
}

# Turn the used range into a proper Excel Table ("Table1") with an
# autofilter on the header row.
$tableRange = $ws.Range("A1:U60")
$tbl = $ws.ListObjects.Add(1, $tableRange, [System.Type]::Missing, 1)
$tbl.Name = "Table1"

# Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
